# Second Problem -- Part B Define
# Defined and broke down the ratios of part b's problem
#
# The original last paragraph ("For the smallest amount of socks to draw a
# matching pair, you'll need only 10 socks.") is turned into a bold topic
# sentence, and six new paragraphs are inserted after it that define and
# break down the sock-color ratios. The original paragraph's run content is
# replaced with the concluding "2/5ths of the socks are NOT black." sentence
# (keeping its original sz/szCs paragraph-mark formatting and the _GoBack
# bookmark), which becomes the new final paragraph in the document.

$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Insert-ParaAfter($afterPara, [string]$innerXml) {
    $d = $word.ActiveDocument
    $r = $afterPara.Range.Duplicate
    $r.InsertParagraphAfter() | Out-Null
    $idx = $afterPara.Index + 1
    $newPara = $d.Paragraphs.Item($idx)
    $nr = $newPara.Range.Duplicate
    $nr.Collapse(0)
    $xml = '<w:p xmlns:w="' + $W + '">' + $innerXml + '</w:p>'
    $nr.InsertXML($xml)
    return $d.Paragraphs.Item($idx)
}

function Set-ParaXml($para, [string]$innerXml) {
    $r = $para.Range.Duplicate
    $r.Collapse(1)
    $xml = '<w:p xmlns:w="' + $W + '">' + $innerXml + '</w:p>'
    $r.InsertXML($xml)
}

$d = $word.ActiveDocument

# Locate the target paragraph ("For the smallest amount of socks...").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*smallest amount of socks*") {
        $target = $cand
        break
    }
}

# --- New paragraph 2: "Define." + ratio explanation --------------------
$p2 = '<w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr>' `
    + '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Define. </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">As for obtaining one matching pair of each color: Since there is a total of 20 </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>socks</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> in a drawer, the chances of obtaining a certain sock color are as follows:</w:t></w:r>'

# --- New paragraph 3: Black bullet --------------------------------------
$p3 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr>' `
    + '<w:r><w:t>Black' + [string][char]0x2014 + '10:20</w:t></w:r>'

# --- New paragraph 4: Brown bullet --------------------------------------
$p4 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr>' `
    + '<w:r><w:t>Brown ' + [string][char]0x2013 + ' 6:20</w:t></w:r>'

# --- New paragraph 5: White bullet --------------------------------------
$p5 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr>' `
    + '<w:r><w:t>White ' + [string][char]0x2013 + ' 4:20</w:t></w:r>'

# --- New paragraph 6: "Thus, if we put..." ------------------------------
$p6 = '<w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr>' `
    + '<w:r><w:t>Thus, if we put this situation in terms of fractions, the following facts are true:</w:t></w:r>'

# --- New paragraph 7: "Half the socks..." bullet ------------------------
$p7 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr></w:pPr>' `
    + '<w:r><w:t>Half the socks in the drawer are black.</w:t></w:r>'

# --- New final paragraph: "2/5ths..." (keeps sz/szCs + _GoBack bookmark) -
$p8 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' `
    + '<w:r><w:t>2/5ths of the socks are NOT black.</w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

# Insert the six new paragraphs after the target paragraph, in order.
$cursor = Insert-ParaAfter $target $p2
$cursor = Insert-ParaAfter $cursor $p3
$cursor = Insert-ParaAfter $cursor $p4
$cursor = Insert-ParaAfter $cursor $p5
$cursor = Insert-ParaAfter $cursor $p6
$cursor = Insert-ParaAfter $cursor $p7
Insert-ParaAfter $cursor $p8 | Out-Null

# Turn the original target paragraph into the bold topic sentence (drop the
# sz/szCs paragraph-mark formatting in favor of bold, and drop the bookmark
# which now lives on the new final paragraph).
$rightQuote = [string][char]0x2019
$p1 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:b/></w:rPr></w:pPr>' `
    + '<w:r><w:t xml:space="preserve">For the smallest amount of socks to draw a matching pair, </w:t></w:r>' `
    + '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">you' + $rightQuote + 'll need </w:t></w:r>' `
    + '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">only </w:t></w:r>' `
    + '<w:r><w:rPr><w:b/></w:rPr><w:t>10</w:t></w:r>' `
    + '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> socks.</w:t></w:r>'

Set-ParaXml $target $p1
